$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 13.17322371252606

# Row 3
$ws.Range("F3").Value = 8.25975769155853

# Row 4
$ws.Range("B4").Value = 27
$ws.Range("C4").Value = "5ff8ad350d084e10f500e48a"
$ws.Range("D4").Value = "Drew"
$ws.Range("F4").Value = 7.198709993617562

# Row 5
$ws.Range("B5").Value = 30
$ws.Range("C5").Value = "60c2341fe95d71ee52c043f0"
$ws.Range("D5").Value = "Matthew"
$ws.Range("F5").Value = 7.013420770724821

# Row 6
$ws.Range("F6").Value = 5.30118687809812

# Row 7
$ws.Range("F7").Value = 5.217927984380697

# Row 8
$ws.Range("F8").Value = 5.022934074744907

# Row 9
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = "5e2522d6b734b47915f88275"
$ws.Range("D9").Value = "Corey"
$ws.Range("F9").Value = 4.346124827430741

# Row 10
$ws.Range("B10").Value = 33
$ws.Range("C10").Value = "60b322994d0b901954690036"
$ws.Range("D10").Value = "Brennan"
$ws.Range("F10").Value = 4.218944548501164

# Row 11
$ws.Range("F11").Value = 3.102299781206951

# Row 12
$ws.Range("B12").Value = 29
$ws.Range("C12").Value = "60b83826821417f8e484a207"
$ws.Range("D12").Value = "Eli"
$ws.Range("F12").Value = 2.368353172506046
$ws.Range("G12").Value = "White"

# Row 13
$ws.Range("B13").Value = 50
$ws.Range("C13").Value = "6097b95056caf5ebb2720002"
$ws.Range("D13").Value = "Damian"
$ws.Range("F13").Value = 2.301699192143767
$ws.Range("G13").Value = "Black or African American"

